$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Kitchen: Number of burners? (R2) was stored as text "4"; change it to the
# numeric value 3.
$ws.Range("R2").Value = 3

# Update the view: zoom to 140% and move the active selection to R3.
$ws.Application.ActiveWindow.Zoom = 140
$ws.Range("R3").Select()
